# Optuna Attempt (go back with original)
# Updates forecast values on "Forecast Comparison" and the derived totals on
# "Summary" to reflect the re-run forecast numbers.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# MyForecast (D), Inventory Coverage (H), Seasonality Index (L)
$ws1.Range("L2").Value  = 1.2

$ws1.Range("D3").Value  = 204
$ws1.Range("H3").Value  = 4.46
$ws1.Range("L3").Value  = 0.89

$ws1.Range("H4").Value  = 4.56
$ws1.Range("L4").Value  = 1.06

$ws1.Range("D5").Value  = 148
$ws1.Range("H5").Value  = 3.72
$ws1.Range("L5").Value  = 1.19

$ws1.Range("D6").Value  = 145
$ws1.Range("H6").Value  = 2.78
$ws1.Range("L6").Value  = 1.17

$ws1.Range("D7").Value  = 144
$ws1.Range("H7").Value  = 1.79
$ws1.Range("L7").Value  = 0.98

$ws1.Range("D8").Value  = 123
$ws1.Range("H8").Value  = 0.93
$ws1.Range("L8").Value  = 0.86

$ws1.Range("D9").Value  = 88
$ws1.Range("L9").Value  = 1.18

$ws1.Range("D10").Value = 141
$ws1.Range("L10").Value = 1.2

$ws1.Range("D11").Value = 132
$ws1.Range("L11").Value = 1.19

$ws1.Range("D12").Value = 72
$ws1.Range("L12").Value = 1.13

$ws1.Range("D13").Value = 132
$ws1.Range("L13").Value = 0.93

$ws1.Range("D14").Value = 128
$ws1.Range("L14").Value = 0.93

$ws1.Range("D15").Value = 123
$ws1.Range("L15").Value = 0.92

$ws1.Range("D16").Value = 124
$ws1.Range("L16").Value = 1.08

$ws1.Range("D17").Value = 123
$ws1.Range("L17").Value = 1.17

# Summary totals are stored as text (numbers-as-strings), so force the
# "@" text format before writing and drop back to the Normal style
# afterwards so no stray number format lingers on the cell.
$ws2 = $wb.Worksheets.Item("Summary")

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = [string]$value
    $range.Style = "Normal"
}

Set-TextValue $ws2.Range("B9")  2175
Set-TextValue $ws2.Range("B10") 1195
Set-TextValue $ws2.Range("B11") 695
Set-TextValue $ws2.Range("B12") 204
Set-TextValue $ws2.Range("B14") 72
